# Apply the commit's changes:
#  - Rename "Sheet1" to "Blockers"
#  - Add a new "Observations" sheet after "Blockers" with header + one data row
#  - Re-use the existing header styling (bold) and date number format
#  - Give the "Transaction Details" cell in column D its own (black) font, matching
#    the new font/style introduced in the edit
#  - Adjust the Observations sheet column E width
#  - Leave the new sheet active/selected, matching the saved view state

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Blockers"

# Add the new sheet right after "Blockers"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Observations"

# ---- Header row (row 1) ----
$ws2.Range("A1").Value = "S.no"
$ws2.Range("B1").Value = "Date "
$ws2.Range("C1").Value = "Feature"
$ws2.Range("D1").Value = "Sub Feature"
$ws2.Range("E1").Value = "Element Name "
$ws2.Range("F1").Value = "Effected Functionality"
$ws2.Range("A1:F1").Font.Bold = $true

# ---- Data row (row 2) ----
$ws2.Range("A2").Value = 1

$ws2.Range("B2").Value = 45132
$ws2.Range("B2").NumberFormat = $ws1.Range("B2").NumberFormat

$ws2.Range("C2").Value = "Transaction Details"

$ws2.Range("D2").Value = "Transaction Details"
$ws2.Range("D2").Font.Color = 0

$ws2.Range("E2").Value = "Buy Tokens-Debit,Credit"
$ws2.Range("F2").Value = "its showing Buy Token instead of Buy Tokens in Debit,Credit Transaction details"

# ---- Column widths ----
$ws2.Columns.Item(5).ColumnWidth = 26.5
$ws1.Columns.Item(6).ColumnWidth = 21.999

# ---- View state: Observations becomes the active/visible tab ----
$ws2.Activate()
$ws2.Range("E2").Select()
